$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "307.55"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.10%"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.11"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.15%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.120"
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.17%"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07797"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.23%"
$ws.Range("E5").Style = "Normal"

$ws.Range("B6").Value = "GateToken"

$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.421"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.04%"
$ws.Range("E6").Style = "Normal"

$ws.Range("B7").Value = "KuCoinToken"

$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.259"
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.50%"
$ws.Range("E7").Style = "Normal"

$ws.Range("B8").Value = "FTXToken"

$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.887"
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.05%"
$ws.Range("E8").Style = "Normal"

$ws.Range("B9").Value = "BTSEToken"

$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.991"
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "8.73%"
$ws.Range("E9").Style = "Normal"

$ws.Range("B10").Value = "MXToken"

$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9254"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.32%"
$ws.Range("E10").Style = "Normal"

$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"

$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1083"
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-10.92%"
$ws.Range("E11").Style = "Normal"

$ws.Range("B12").Value = "WazirX"

$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1899"
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.63%"
$ws.Range("E12").Style = "Normal"

$ws.Range("B13").Value = "MandalaExchangeToken"

$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08890"
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-3.43%"
$ws.Range("E13").Style = "Normal"

$ws.Range("B14").Value = "BitrueCoin"

$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03328"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-2.98%"
$ws.Range("E14").Style = "Normal"

$ws.Range("B15").Value = "BitMartToken"

$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09579"
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.02%"
$ws.Range("E15").Style = "Normal"

$ws.Range("B16").Value = "BitForexToken"

$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001382"
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.01%"
$ws.Range("E16").Style = "Normal"

$ws.Range("B17").Value = "TigerCash"

$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005799"
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.71%"
$ws.Range("E17").Style = "Normal"

$ws.Range("B18").Value = "LEO"

$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.510"
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.39%"
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3416"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.34%"
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.292"
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "19.70%"
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1276"
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.45%"
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2504"
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-3.40%"
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04382"
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.46%"
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.45%"
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004252"
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.11%"
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.09%"
$ws.Range("E26").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02155"
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.10%"
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04999"
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.27%"
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007477"
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.60%"
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1347"
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.20%"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008651"
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-11.39%"
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.99%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008000"
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006535"
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.79%"
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.08%"
$ws.Range("E47").Style = "Normal"

$ws.Range("B48").Value = "BOLO"

$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002854"
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-2.84%"
$ws.Range("E48").Style = "Normal"

$ws.Range("B49").Value = "CoinbaseStockToken"

$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001002"
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-16.60%"
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.08%"
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.08%"
$ws.Range("E51").Style = "Normal"
